$d = $word.ActiveDocument

$replacements = @(
    @{old = "2026-01-05 Monday"; new = "2026-01-06 Tuesday"},
    @{old = "289÷2="; new = "802÷9="},
    @{old = "110÷6="; new = "952÷7="},
    @{old = "838÷7="; new = "269÷4="},
    @{old = "721÷9="; new = "822÷3="},
    @{old = "464÷8="; new = "348÷8="},
    @{old = "742÷6="; new = "485÷6="},
    @{old = "967÷3="; new = "931÷5="},
    @{old = "793÷9="; new = "103÷7="},
    @{old = "627÷3="; new = "583÷7="},
    @{old = "145÷6="; new = "452÷5="},
    @{old = "981÷8="; new = "176÷2="},
    @{old = "200÷9="; new = "107÷3="},
    @{old = "992÷9="; new = "997÷2="},
    @{old = "612÷7="; new = "437÷2="},
    @{old = "684÷5="; new = "678÷8="},
    @{old = "244÷4="; new = "357÷4="},
    @{old = "468÷5="; new = "464÷2="},
    @{old = "914÷7="; new = "965÷5="},
    @{old = "812÷3="; new = "604÷8="},
    @{old = "101÷4="; new = "406÷7="},
    @{old = "960÷3="; new = "580÷8="},
    @{old = "192÷4="; new = "667÷2="},
    @{old = "419÷7="; new = "874÷9="},
    @{old = "788÷8="; new = "565÷2="},
    @{old = "468÷9="; new = "808÷2="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
